$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Post Test" (column C) values for rows 2-23
$postTestValues = @(0, 0, 0, 1, 1, 3, 2, 3, 1, 2, 2, 1, 0, 0, 1, 0, 0, 1, 0, 3, 2, 1)

for ($i = 0; $i -lt $postTestValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $postTestValues[$i]
}

# Update the active selection to C24, matching the sheetView selection change
$ws.Range("C24").Select()
